# Regenerate the lattice-multiplication exercise table: every cell's
# problem (N x M), the two factor digits, the divider, and the two
# lattice-row labels are replaced with a freshly generated exercise,
# while keeping the existing 5 row x 3 column grid and run formatting
# (sz 32) untouched.

function Set-LatticeCell {
    param($Table, $Row, $Col, $Values)
    $vt = [char]11
    $joined = [string]::Join($vt, $Values)
    $cell = $Table.Cell($Row, $Col)
    $cell.Range.Text = $joined
}

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

Set-LatticeCell $t 1 1 @("32 x 88", "  8    8", "  ----", "3|    |", "2|    |")
Set-LatticeCell $t 1 2 @("43 x 54", "  5    4", "  ----", "4|    |", "3|    |")
Set-LatticeCell $t 1 3 @("98 x 68", "  6    8", "  ----", "9|    |", "8|    |")

Set-LatticeCell $t 2 1 @("38 x 28", "  2    8", "  ----", "3|    |", "8|    |")
Set-LatticeCell $t 2 2 @("48 x 97", "  9    7", "  ----", "4|    |", "8|    |")
Set-LatticeCell $t 2 3 @("68 x 45", "  4    5", "  ----", "6|    |", "8|    |")

Set-LatticeCell $t 3 1 @("39 x 61", "  6    1", "  ----", "3|    |", "9|    |")
Set-LatticeCell $t 3 2 @("44 x 68", "  6    8", "  ----", "4|    |", "4|    |")
Set-LatticeCell $t 3 3 @("88 x 58", "  5    8", "  ----", "8|    |", "8|    |")

Set-LatticeCell $t 4 1 @("31 x 36", "  3    6", "  ----", "3|    |", "1|    |")
Set-LatticeCell $t 4 2 @("63 x 46", "  4    6", "  ----", "6|    |", "3|    |")
Set-LatticeCell $t 4 3 @("45 x 98", "  9    8", "  ----", "4|    |", "5|    |")

Set-LatticeCell $t 5 1 @("37 x 31", "  3    1", "  ----", "3|    |", "7|    |")
Set-LatticeCell $t 5 2 @("98 x 40", "  4    0", "  ----", "9|    |", "8|    |")
Set-LatticeCell $t 5 3 @("37 x 86", "  8    6", "  ----", "3|    |", "7|    |")

Write-Output "Lattice multiplication table regenerated"
